$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Main")
$wsModel = $wb.Worksheets.Item("Model")

# --- Main sheet: update Price (K2) input from 251 to 310 ---
$wsMain.Range("K2").Value = 310

# --- Model sheet: update terminal growth rate ("Maturity", X20) from 1% to 0% ---
$wsModel.Range("X20").Value = 0

# --- Update selections to match the saved view state ---
# Model sheet: bottom-right frozen pane selection becomes S24:V26 (active cell S24)
$wsModel.Range("S24:V26").Select() | Out-Null

# Main sheet: selection becomes K3, and Main should remain the active/selected tab
$wsMain.Range("K3").Select() | Out-Null
